$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Existing "step 2" contents for the three affected test cases (TC2, TC3, TC4).
$filterDesc   = "Chefe Seleciona um usuário para filtrar as autorizações de pagamento associadas a ele; e Submete a busca ao sistema."
$filterResult = "SYSTEM Filtra os registros (autorizações de pagamento pendentes) e exibe apenas aqueles atribuídos ao usuário selecionado."

$authDesc     = "Chefe Clica para realizar a autorização de pagamento."
$authResult   = "SYSTEM Apresenta a tela de Registrar Autorizações de Pagamento"

$assignDesc   = "Chefe Dado um registro selecionado (solicitação aguardando autorização de pagamento - AP), o usuário pode atribuir/desatribuir a responsabilidade da AP a si próprio; e Clica para atribuir/desatribuir o registro a si mesmo."
$assignResult = "SYSTEM Atualiza a lista de registros de solicitações, onde o nome deverá constar o nome do usuário logado (que se atribuiu como responsável pela AP) no campo de atribuição (no caso de desatribuição, o nome deverá ser removido)."

# Rotate the "step 2" row content between the TC2, TC3 and TC4 blocks:
#   TC2 (row 20) gets the assign/unassign content (previously under TC4)
#   TC3 (row 28) gets the filter content (previously under TC2)
#   TC4 (row 36) gets the authorization content (previously under TC3)
$ws.Range("B20").Value = $assignDesc
$ws.Range("D20").Value = $assignResult

$ws.Range("B28").Value = $filterDesc
$ws.Range("D28").Value = $filterResult

$ws.Range("B36").Value = $authDesc
$ws.Range("D36").Value = $authResult
